$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.323.53"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").Value = "3.311.29"
$ws.Range("E3").Value = "  -2.95%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.10%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -3.02%  "
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("E11").Value = "  -3.23%  "
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").Value = "3.881.09"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "3.306.09"
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("E17").Value = "  -2.83%  "
$ws.Range("D18").Value = "60.313.83"
$ws.Range("E18").Value = "  -2.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.44%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("E24").Value = "  -3.76%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -3.74%  "
$ws.Range("E27").Value = "  -7.51%  "
$ws.Range("E28").Value = "  -4.49%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -5.20%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  -2.77%  "
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("E34").Value = "  -2.25%  "
$ws.Range("E35").Value = "  -5.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.24%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.58%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.13%  "
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.339.02"
$ws.Range("E41").Value = "  -3.00%  "
$ws.Range("E42").Value = "  -5.49%  "
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E44").Value = "  -4.10%  "
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("E46").Value = "  -4.46%  "
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("D48").Value = "2.372.59"
$ws.Range("E48").Value = "  -6.61%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E50").Value = "  -5.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.88%  "
